$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $range = $d.Content
    $found = $range.Find.Execute($old, $true, $true, $false, $false, $false, $true, 1, $false, $new, 2)
    if (-not $found) {
        Write-Host "NOT FOUND: $old"
    }
}

Replace-Text "2024-06-29 Saturday" "2024-06-30 Sunday"
Replace-Text "63×15=" "96×61="
Replace-Text "61×61=" "29×21="
Replace-Text "60×75=" "57×82="
Replace-Text "23×93=" "38×18="
Replace-Text "62×88=" "48×48="
Replace-Text "49×99=" "90×41="
Replace-Text "26×12=" "48×93="
Replace-Text "99×73=" "97×73="
Replace-Text "63×33=" "83×63="
Replace-Text "34×75=" "72×93="
Replace-Text "34×91=" "36×18="
Replace-Text "20×18=" "26×98="
Replace-Text "69×54=" "37×96="
Replace-Text "91×77=" "19×46="
Replace-Text "37×38=" "42×90="
Replace-Text "24×31=" "25×71="
Replace-Text "25×81=" "80×23="
Replace-Text "54×72=" "32×88="
Replace-Text "34×76=" "67×25="
Replace-Text "66×55=" "36×56="
Replace-Text "41×71=" "80×79="
Replace-Text "66×11=" "50×27="
Replace-Text "35×35=" "74×85="
Replace-Text "61×62=" "79×51="
Replace-Text "36×13=" "21×45="

Write-Host "Done"
